$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new record at row 20 (date 2021-09-13 / serial 44452),
# pushing all subsequent rows down by one.
$ws.Rows("20:20").Insert()

$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44452
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 300000001
$ws.Range("G20").Value = "Rabanito"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 7000
$ws.Range("N20").Value = "$/docena de paquetes"
$ws.Range("O20").Value = "Provincia de Cautín"
$ws.Range("P20").Value = 583
$ws.Range("Q20").Value = 12
$ws.Range("R20").Value = "Hortaliza"

# Match the number format of the other date cells in column D (style index 2).
$ws.Range("D20").NumberFormat = $ws.Range("D21").NumberFormat
